# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (becomes the 2nd
#    tab), pushing every other quarter sheet down by one position.
# 2. Populate it with the 2022-Q3 fund holdings table.
# 3. Prepend a "2022-Q3" row to the "总计" summary sheet (and keep the
#    existing rows, shifting the trailing "2020-Q4" row from r8 -> r9).
# 4. Restore the originally-active sheet so workbook view state is
#    untouched by the structural edit.

$wb = $excel.ActiveWorkbook

# Remember which sheet was active before we start shuffling tabs around,
# so we can restore that at the end (Worksheets.Add activates the new sheet).
$origActive = $wb.ActiveSheet.Name

$total = $wb.Worksheets.Item(1)           # "总计"
$q2Sheet = $wb.Worksheets.Item(2)         # currently "2022-Q2" (will become 3rd tab)

# --- 1. Insert the new sheet right before the current 2nd tab ---------
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Borrow the header / index-column formatting from the existing,
# identically-laid-out "2022-Q2" sheet (A1:H16 range shape: header row +
# 15 fund rows) so fonts/borders/alignment match the sibling quarter
# sheets exactly.
$q2Sheet.Range("A1:H2").Copy()
$q3Sheet.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$q3Sheet.Range("A2").Copy()
$q3Sheet.Range("A3:A16").PasteSpecial(-4122)
$q3Sheet.Range("H2").Copy()
$q3Sheet.Range("H3:H16").PasteSpecial(-4122)

# --- 2. Header row -------------------------------------------------------
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Columns B-G carry text-shaped values (fund codes, percentages, etc. are
# stored as literal text, not numbers) in every quarter sheet - force the
# "Text" number format before writing so "011578" / "21.81" / "94.20" etc.
# round-trip byte-for-byte instead of being auto-coerced into numbers.
$q3Sheet.Range("B2:G16").NumberFormat = "@"

$q3Data = @(
    @(0,  "011578", "汇丰晋信核心成长混合A", "21.81", "94.20", "3.55", "0.7743", 10),
    @(1,  "210009", "金鹰核心资源混合", "2.84", "93.42", "6.31", "0.1792", 1),
    @(2,  "162102", "金鹰中小盘精选混合", "3.17", "78.28", "5.45", "0.1728", 1),
    @(3,  "001167", "金鹰科技创新股票", "2.66", "94.84", "6.33", "0.1684", 1),
    @(4,  "011579", "汇丰晋信核心成长混合C", "4.19", "94.20", "3.55", "0.1487", 10),
    @(5,  "000458", "英大领先回报混合", "1.32", "93.57", "3.35", "0.0442", 2),
    @(6,  "012200", "新华鑫科技3个月滚动持有灵活配置混合A", "1.28", "71.16", "2.93", "0.0375", 9),
    @(7,  "004890", "中邮健康文娱灵活配置混合", "0.41", "86.15", "3.93", "0.0161", 7),
    @(8,  "012201", "新华鑫科技3个月滚动持有灵活配置混合C", "0.45", "71.16", "2.93", "0.0132", 9),
    @(9,  "001270", "英大灵活配置混合A", "0.28", "93.98", "3.33", "0.0093", 2),
    @(10, "001271", "英大灵活配置混合B", "0.28", "93.98", "3.33", "0.0093", 2),
    @(11, "165524", "信诚中证智能家居指数（LOF）A", "0.35", "91.70", "1.10", "0.0038", 10),
    @(12, "001914", "中信建投聚利混合A", "0.10", "39.73", "2.01", "0.0020", 9),
    @(13, "013084", "信诚中证智能家居指数（LOF）C", "0.15", "91.70", "1.10", "0.0016", 10),
    @(14, "006845", "中信建投聚利混合C", "0.01", "39.73", "2.01", "0.0002", 9)
)

$r = 2
foreach ($row in $q3Data) {
    $q3Sheet.Range("A$r").Value = $row[0]
    $q3Sheet.Range("B$r").Value = $row[1]
    $q3Sheet.Range("C$r").Value = $row[2]
    $q3Sheet.Range("D$r").Value = $row[3]
    $q3Sheet.Range("E$r").Value = $row[4]
    $q3Sheet.Range("F$r").Value = $row[5]
    $q3Sheet.Range("G$r").Value = $row[6]
    $q3Sheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# --- 3. Update the "总计" roll-up sheet ----------------------------------
# Push the "2020-Q4" row from r8 down to r9 (copy formats first so the
# A9 index cell keeps the bordered/bold look of A2:A8), then overwrite
# rows 2-8 top-down with the shifted data and drop the new 2022-Q3 row
# in at the top.
$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)

$total.Range("A9").Value = 7
$total.Range("B9").Value = "2020-Q4"
$total.Range("C9").Value = 4
$total.Range("D9").Value = 0.05

$total.Range("A8").Value = 6
$total.Range("B8").Value = "2021-Q1"
$total.Range("C8").Value = 11
$total.Range("D8").Value = 0.64

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 37
$total.Range("D7").Value = 13.86

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.02

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 21
$total.Range("D5").Value = 7.92

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 6
$total.Range("D4").Value = 0.62

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 13
$total.Range("D3").Value = 0.89

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 15
$total.Range("D2").Value = 1.58

# --- 4. Restore the view state (active tab) ------------------------------
$wb.Worksheets.Item($origActive).Activate()
